# Weekly update: insert this week's new price record at the top of the data
# (row 2), pushing all previously-recorded rows down by one. The worksheet
# only has a single data table starting at row 2 (row 1 is the header).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row right above the current first data row (row 2).
# This shifts every existing data row (old rows 2-35) down to rows 3-36,
# and Excel automatically extends the sheet dimension / used range to
# A1:R36.
$ws.Rows.Item(2).Insert()

# The newly inserted row inherits formatting from the row above it (the
# bold header row). Reset it to the plain/default formatting used by the
# rest of the data rows.
$ws.Range("A2:R2").ClearFormats()

# Column D holds dates stored as date-formatted numbers; restore that
# number format for the new row's date cell to match the other rows.
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Populate the new row with this week's record.
$arr = New-Object 'object[,]' 1,18
$arr[0,0]  = 11
$arr[0,1]  = 'Vega Monumental Concepción'
$arr[0,2]  = 'Bíobío'
$arr[0,3]  = 44631
$arr[0,4]  = 8
$arr[0,5]  = 100112030
$arr[0,6]  = 'Poroto granado'
$arr[0,7]  = 'Sin especificar'
$arr[0,8]  = 'Primera'
$arr[0,9]  = 110
$arr[0,10] = 20000
$arr[0,11] = 21000
$arr[0,12] = 20455
$arr[0,13] = '$/saco 25 kilos'
$arr[0,14] = "Región de O'Higgins"
$arr[0,15] = 818
$arr[0,16] = 25
$arr[0,17] = 'Hortaliza'

$ws.Range("A2:R2").Value = $arr
